# Generate Report for Handback
# A new handback ("93b72175-5d84-432b-ad7a-04e0a288af30") was processed.
# It is inserted as the new row 3 on every sheet (Overview, zh-cn, de-de),
# pushing the previously-row-3 entry ("f04efb91-24d9-4f68-89a4-3b0ba9fc450c")
# down to a newly appended row 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (7 columns: A..G)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

# Row 3 becomes the new handback entry.
$wsOverview.Range("A3").Value = "93b72175-5d84-432b-ad7a-04e0a288af30.md"
$wsOverview.Range("B3").Value = "e2e\93b72175-5d84-432b-ad7a-04e0a288af30.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-25 02:44:07"

$hl = $wsOverview.Range("B3").Hyperlinks.Item(1)
$hl.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/93b72175-5d84-432b-ad7a-04e0a288af30.md"
$hl.TextToDisplay = "e2e\93b72175-5d84-432b-ad7a-04e0a288af30.md"

# Row 4 (new) takes what used to live in row 3.
$wsOverview.Range("A4").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
$wsOverview.Range("B4").Value = "e2e\f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-25 02:39:33"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be83e81ba0665194049ffb60eaf7f18c025090e2/e2e/f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md", "", "", "e2e\f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheets "zh-cn" and "de-de" (16 columns: A..P) share the same layout.
# ---------------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Ext = "zh-cn.xlf" },
    @{ Name = "de-de"; Ext = "de-de.xlf" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null

    # Row 3 becomes the new handback entry (93b72175...).
    $ws.Range("A3").Value = "93b72175-5d84-432b-ad7a-04e0a288af30.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
    $ws.Range("D3").Value = "e2e"
    $ws.Range("E3").Value = "ht"
    $ws.Range("F3").Value = "'True"
    $ws.Range("G3").Value = "93b72175-5d84-432b-ad7a-04e0a288af30.5e12453ae228e8add7253157f710373643a3a46d.$($lang.Ext)"
    if ($lang.Name -eq "zh-cn") {
        $ws.Range("H3").Value = "2016-08-25 02:43:58"
    } else {
        $ws.Range("H3").Value = "2016-08-25 02:44:07"
    }
    $ws.Range("I3").Value = "93b72175-5d84-432b-ad7a-04e0a288af30.md"
    $ws.Range("J3").Value = "93b72175-5d84-432b-ad7a-04e0a288af30.5e12453ae228e8add7253157f710373643a3a46d.$($lang.Ext)"
    if ($lang.Name -eq "zh-cn") {
        $ws.Range("K3").Value = "2016-08-25 02:44:28"
    } else {
        $ws.Range("K3").Value = "2016-08-25 02:44:35"
    }
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = "'True"
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = "'False"
    $ws.Range("P3").Value = ""

    $hlA = $ws.Range("A3").Hyperlinks.Item(1)
    $hlA.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/93b72175-5d84-432b-ad7a-04e0a288af30.md"
    $hlA.TextToDisplay = "93b72175-5d84-432b-ad7a-04e0a288af30.md"

    $hlI = $ws.Range("I3").Hyperlinks.Item(1)
    $hlI.Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/93b72175-5d84-432b-ad7a-04e0a288af30.md"
    $hlI.TextToDisplay = "93b72175-5d84-432b-ad7a-04e0a288af30.md"

    # Row 4 (new) takes what used to live in row 3 (f04efb91...).
    $ws.Range("A4").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Handed back: in sync with en-US"
    $ws.Range("D4").Value = "e2e"
    $ws.Range("E4").Value = "ht"
    $ws.Range("F4").Value = "'True"
    $ws.Range("G4").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.0e0ee87858c68783cd4d11057245d5d9d0c48721.$($lang.Ext)"
    if ($lang.Name -eq "zh-cn") {
        $ws.Range("H4").Value = "2016-08-25 02:39:28"
    } else {
        $ws.Range("H4").Value = "2016-08-25 02:39:33"
    }
    $ws.Range("I4").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
    $ws.Range("J4").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.0e0ee87858c68783cd4d11057245d5d9d0c48721.$($lang.Ext)"
    if ($lang.Name -eq "zh-cn") {
        $ws.Range("K4").Value = "2016-08-25 02:39:45"
    } else {
        $ws.Range("K4").Value = "2016-08-25 02:39:52"
    }
    $ws.Range("L4").Value = ""
    $ws.Range("M4").Value = "'True"
    $ws.Range("N4").Value = ""
    $ws.Range("O4").Value = "'False"
    $ws.Range("P4").Value = ""

    $urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/be83e81ba0665194049ffb60eaf7f18c025090e2/e2e/f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
    $ws.Hyperlinks.Add($ws.Range("A4"), $urlBase, "", "", "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I4"), $urlBase, "", "", "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md") | Out-Null
}
